# Update scripts with new TPM values.
# The underlying data table now has only 3 sending-cluster/target-cluster
# combinations (target cluster fixed to "FAPs"), so rows 5-7 are removed
# and rows 2-4 get the recomputed NATMI metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (old rows 5, 6 and 7).
$ws.Rows.Item(5).Resize(3).Delete()

# Row 2: ECs -> FAPs (Leap2 -> Ghsr)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Leap2"
$ws.Range("C2").Value = "Ghsr"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9420959999999999
$ws.Range("H2").Value = 2.826288
$ws.Range("I2").Value = 0.2074065050464874
$ws.Range("J2").Value = 0.2074065050464874
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1546876666666667
$ws.Range("N2").Value = 0.464063
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.145730632016
$ws.Range("R2").Value = 1.311575688144
$ws.Range("S2").Value = 0.2074065050464874
$ws.Range("T2").Value = 0.2074065050464874

# Row 3: FAPs -> FAPs (Leap2 -> Ghsr)
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Leap2"
$ws.Range("C3").Value = "Ghsr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.587073666666666
$ws.Range("H3").Value = 7.761221
$ws.Range("I3").Value = 0.5695554460491655
$ws.Range("J3").Value = 0.5695554460491655
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1546876666666667
$ws.Range("N3").Value = 0.464063
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.4001883889914444
$ws.Range("R3").Value = 3.601695500923
$ws.Range("S3").Value = 0.5695554460491655
$ws.Range("T3").Value = 0.5695554460491655

# Row 4: MuSCs -> FAPs (Leap2 -> Ghsr)
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Leap2"
$ws.Range("C4").Value = "Ghsr"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.013098666666667
$ws.Range("H4").Value = 3.039296
$ws.Range("I4").Value = 0.223038048904347
$ws.Range("J4").Value = 0.223038048904347
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1546876666666667
$ws.Range("N4").Value = 0.464063
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.1567138688497778
$ws.Range("R4").Value = 1.410424819648
$ws.Range("S4").Value = 0.223038048904347
$ws.Range("T4").Value = 0.223038048904347
